# Change access from public to admin of "active posts"
# - add a new row (row 11) recording the change of "添加管理员查看人气话题的功能"
#   by 雷建坤, referencing user story 17
# - bump the "user story编号" (column D) on rows 8-10 from 8 to 16
#   (those 3 existing rows also now belong to story 16)
# - adjust the saved selection to match the authored file

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the user-story numbers for the existing "添加XX管理功能" rows (8-10)
$ws.Range("D8").Value = 16
$ws.Range("D9").Value = 16
$ws.Range("D10").Value = 16

# Append the new entry for admin-only access to popular/active topics.
# Insert the new row by copying row 10's formatting (keeps the same date
# style as the rest of the "日期" column) and then overwrite its values.
$ws.Rows("10:10").Copy()
$ws.Rows("11:11").Insert()
$ws.Range("A11").Value = 10
$ws.Range("C11").Value = "雷建坤"
$ws.Range("D11").Value = 17
$ws.Range("E11").Value = "添加管理员查看人气话题的功能"

# Match the saved selection from the authored workbook
$ws.Range("E15").Select()
